# Auto commit at 2025-10-13  7:47:43.56
#
# Updates the "Metrics" sheet's daily figures (B2:B13) with refreshed
# totals, which ripple via formulas into the "today" sheet (B11:B22,
# E11:E22, F11:F22) and the TODAY()-1 driven date in A1. Also moves the
# active-sheet/selection bookmarks: "today" had been the active tab with
# D4 selected; "Metrics" becomes active with D17 selected, and "today"'s
# own remembered selection moves to J12.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

# --- refresh the Metrics figures -------------------------------------
$metrics.Range("B2").Value  = 182490.06000000003
$metrics.Range("B3").Value  = 152583.63
$metrics.Range("B4").Value  = 65365.56
$metrics.Range("B5").Value  = 7141
$metrics.Range("B6").Value  = 4549621.53
$metrics.Range("B7").Value  = 3842402.3
$metrics.Range("B8").Value  = 1335967.7000000002
$metrics.Range("B9").Value  = 176142
$metrics.Range("B10").Value = 33014945.330999829
$metrics.Range("B11").Value = 31117623.820000004
$metrics.Range("B12").Value = 11617676.59
$metrics.Range("B13").Value = 1273769

# --- move "today"'s remembered selection to J12, then hand the active
#     tab / tabSelected flag over to "Metrics" at D17 ------------------
$today.Select()
$today.Range("J12").Select()

$metrics.Select()
$metrics.Range("D17").Select()
